# "add new words in excel"
# Adds a new daily worksheet ("2025-08-11") right after the existing
# "2025-08-10" sheet, fills it with that day's vocabulary list (Chinese
# term in column A, English term in column B), and makes the new sheet
# the active/selected one (mirroring the author switching tabs to the
# freshly-created sheet before saving).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New worksheet, inserted immediately after "2025-08-10"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "2025-08-11"

# Row data: Chinese term -> English term (row 26 only has a column A value)
$ws2.Cells.Item(1, 1).Value  = '有血有肉的'
$ws2.Cells.Item(1, 2).Value  = 'flesh-and-blood'
$ws2.Cells.Item(2, 1).Value  = '撤回'
$ws2.Cells.Item(2, 2).Value  = 'retraction'
$ws2.Cells.Item(3, 1).Value  = '揭露丑闻'
$ws2.Cells.Item(3, 2).Value  = 'muckraking'
$ws2.Cells.Item(4, 1).Value  = '笨重的'
$ws2.Cells.Item(4, 2).Value  = 'bulky'
$ws2.Cells.Item(5, 1).Value  = '预言'
$ws2.Cells.Item(5, 2).Value  = 'divine'
$ws2.Cells.Item(6, 1).Value  = '追踪'
$ws2.Cells.Item(6, 2).Value  = 'trace'
$ws2.Cells.Item(7, 1).Value  = '挖掘'
$ws2.Cells.Item(7, 2).Value  = 'excavations'
$ws2.Cells.Item(8, 1).Value  = '吝啬'
$ws2.Cells.Item(8, 2).Value  = 'parsimony'
$ws2.Cells.Item(9, 1).Value  = '投机主义'
$ws2.Cells.Item(9, 2).Value  = 'opportunism'
$ws2.Cells.Item(10, 1).Value = '吝啬'
$ws2.Cells.Item(10, 2).Value = 'miserliness'
$ws2.Cells.Item(11, 1).Value = '讨厌的东西'
$ws2.Cells.Item(11, 2).Value = 'nuisance'
$ws2.Cells.Item(12, 1).Value = '意外收获'
$ws2.Cells.Item(12, 2).Value = 'windfall'
$ws2.Cells.Item(13, 1).Value = '引起'
$ws2.Cells.Item(13, 2).Value = 'invoke'
$ws2.Cells.Item(14, 1).Value = '同意'
$ws2.Cells.Item(14, 2).Value = 'concur'
$ws2.Cells.Item(15, 1).Value = '井然有序的'
$ws2.Cells.Item(15, 2).Value = 'methodical'
$ws2.Cells.Item(16, 1).Value = '尽管'
$ws2.Cells.Item(16, 2).Value = 'that said'
$ws2.Cells.Item(17, 1).Value = '开玩笑'
$ws2.Cells.Item(17, 2).Value = 'jest'
$ws2.Cells.Item(18, 1).Value = '开放性的'
$ws2.Cells.Item(18, 2).Value = 'expansive'
$ws2.Cells.Item(19, 1).Value = '细小的'
$ws2.Cells.Item(19, 2).Value = 'fine'
$ws2.Cells.Item(20, 1).Value = '尴尬'
$ws2.Cells.Item(20, 2).Value = 'discomfiture'
# NOTE: column B is entered before column A on this row so the shared-string
# table ends up with the same id ordering the workbook was saved with.
$ws2.Cells.Item(21, 2).Value = 'presumptuousness'
$ws2.Cells.Item(21, 1).Value = '自以为是'
$ws2.Cells.Item(22, 1).Value = '细心'
$ws2.Cells.Item(22, 2).Value = 'circumspection'
$ws2.Cells.Item(23, 1).Value = '颁布'
$ws2.Cells.Item(23, 2).Value = 'promulgated'
$ws2.Cells.Item(24, 1).Value = '误解'
$ws2.Cells.Item(24, 2).Value = 'misconstrued'
$ws2.Cells.Item(25, 1).Value = '壁画'
$ws2.Cells.Item(25, 2).Value = 'mural'
$ws2.Cells.Item(26, 1).Value = 's'

# Leave the cursor on the last entered row/cell, like after manual typing
[void]$ws2.Range("A26").Select()

# Switch to the new sheet so it becomes the active/selected tab on save
[void]$ws2.Activate()
